# Fill in the Bible Verse / Lead / Message for the 2025-10-30 study (row 7)
# and leave the selection where the author last left it (B8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = "Nkauj Qhuas Vajtswv 40: 12-13"
$ws.Range("C7").Value = "Sheng"
$ws.Range("D7").Value = "Cia siab rau Vajtswv thaum yus muaj kev nyuab siab"

$ws.Range("B8").Select() | Out-Null
